$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 91.94136433333334
$ws.Range("H2").Value = 275.824093
$ws.Range("I2").Value = 0.02307547609860541
$ws.Range("J2").Value = 0.02307547609860541
$ws.Range("M2").Value = 6.066605666666667
$ws.Range("N2").Value = 18.199817
$ws.Range("O2").Value = 0.849784628791665
$ws.Range("P2").Value = 0.8497846287916652
$ws.Range("Q2").Value = 557.7720018656646
$ws.Range("R2").Value = 5019.948016790981
$ws.Range("S2").Value = 0.01960918489064434
$ws.Range("T2").Value = 0.01960918489064434
$ws.Range("G3").Value = 91.94136433333334
$ws.Range("H3").Value = 275.824093
$ws.Range("I3").Value = 0.02307547609860541
$ws.Range("J3").Value = 0.02307547609860541
$ws.Range("O3").Value = 0.1196497582104962
$ws.Range("P3").Value = 0.1196497582104962
$ws.Range("Q3").Value = 78.53435199776088
$ws.Range("R3").Value = 706.809167979848
$ws.Range("S3").Value = 0.002760975135790221
$ws.Range("T3").Value = 0.002760975135790222
$ws.Range("G4").Value = 91.94136433333334
$ws.Range("H4").Value = 275.824093
$ws.Range("I4").Value = 0.02307547609860541
$ws.Range("J4").Value = 0.02307547609860541
$ws.Range("M4").Value = 0.1824346666666667
$ws.Range("N4").Value = 0.547304
$ws.Range("O4").Value = 0.02555468148257719
$ws.Range("P4").Value = 0.02555468148257719
$ws.Range("Q4").Value = 16.77329215503022
$ws.Range("R4").Value = 150.959629395272
$ws.Range("S4").Value = 0.0005896864417586842
$ws.Range("T4").Value = 0.0005896864417586843
$ws.Range("G5").Value = 91.94136433333334
$ws.Range("H5").Value = 275.824093
$ws.Range("I5").Value = 0.02307547609860541
$ws.Range("J5").Value = 0.02307547609860541
$ws.Range("M5").Value = 0.035773
$ws.Range("N5").Value = 0.107319
$ws.Range("O5").Value = 0.005010931515261538
$ws.Range("P5").Value = 0.005010931515261539
$ws.Range("Q5").Value = 3.289018426296333
$ws.Range("R5").Value = 29.601165836667
$ws.Range("S5").Value = 0.0001156296304121662
$ws.Range("T5").Value = 0.0001156296304121662
$ws.Range("I6").Value = 0.9681738695089209
$ws.Range("J6").Value = 0.9681738695089209
$ws.Range("M6").Value = 6.066605666666667
$ws.Range("N6").Value = 18.199817
$ws.Range("O6").Value = 0.849784628791665
$ws.Range("P6").Value = 0.8497846287916652
$ws.Range("Q6").Value = 23402.34606828563
$ws.Range("R6").Value = 210621.1146145707
$ws.Range("S6").Value = 0.8227392723064283
$ws.Range("T6").Value = 0.8227392723064284
$ws.Range("I7").Value = 0.9681738695089209
$ws.Range("J7").Value = 0.9681738695089209
$ws.Range("O7").Value = 0.1196497582104962
$ws.Range("P7").Value = 0.1196497582104962
$ws.Range("S7").Value = 0.1158417693924629
$ws.Range("T7").Value = 0.1158417693924629
$ws.Range("I8").Value = 0.9681738695089209
$ws.Range("J8").Value = 0.9681738695089209
$ws.Range("M8").Value = 0.1824346666666667
$ws.Range("N8").Value = 0.547304
$ws.Range("O8").Value = 0.02555468148257719
$ws.Range("P8").Value = 0.02555468148257719
$ws.Range("Q8").Value = 703.7541977788568
$ws.Range("R8").Value = 6333.787780009711
$ws.Range("S8").Value = 0.02474137485505472
$ws.Range("T8").Value = 0.02474137485505472
$ws.Range("I9").Value = 0.9681738695089209
$ws.Range("J9").Value = 0.9681738695089209
$ws.Range("M9").Value = 0.035773
$ws.Range("N9").Value = 0.107319
$ws.Range("O9").Value = 0.005010931515261538
$ws.Range("P9").Value = 0.005010931515261539
$ws.Range("Q9").Value = 137.9967929184313
$ws.Range("R9").Value = 1241.971136265882
$ws.Range("S9").Value = 0.004851452954974964
$ws.Range("T9").Value = 0.004851452954974965
$ws.Range("G10").Value = 1.840730666666667
$ws.Range("H10").Value = 5.522192
$ws.Range("I10").Value = 0.000461987232956876
$ws.Range("J10").Value = 0.000461987232956876
$ws.Range("M10").Value = 6.066605666666667
$ws.Range("N10").Value = 18.199817
$ws.Range("O10").Value = 0.849784628791665
$ws.Range("P10").Value = 0.8497846287916652
$ws.Range("Q10").Value = 11.16698709320711
$ws.Range("R10").Value = 100.502883838864
$ws.Range("S10").Value = 0.0003925896492647473
$ws.Range("T10").Value = 0.0003925896492647474
$ws.Range("G11").Value = 1.840730666666667
$ws.Range("H11").Value = 5.522192
$ws.Range("I11").Value = 0.000461987232956876
$ws.Range("J11").Value = 0.000461987232956876
$ws.Range("O11").Value = 0.1196497582104962
$ws.Range("P11").Value = 0.1196497582104962
$ws.Range("Q11").Value = 1.572312866545778
$ws.Range("R11").Value = 14.150815798912
$ws.Range("S11").Value = 0.00005527666071962638
$ws.Range("T11").Value = 0.0000552766607196264
$ws.Range("G12").Value = 1.840730666666667
$ws.Range("H12").Value = 5.522192
$ws.Range("I12").Value = 0.000461987232956876
$ws.Range("J12").Value = 0.000461987232956876
$ws.Range("M12").Value = 0.1824346666666667
$ws.Range("N12").Value = 0.547304
$ws.Range("O12").Value = 0.02555468148257719
$ws.Range("P12").Value = 0.02555468148257719
$ws.Range("Q12").Value = 0.3358130855964445
$ws.Range("R12").Value = 3.022317770368
$ws.Range("S12").Value = 0.00001180593658723015
$ws.Range("T12").Value = 0.00001180593658723015
$ws.Range("G13").Value = 1.840730666666667
$ws.Range("H13").Value = 5.522192
$ws.Range("I13").Value = 0.000461987232956876
$ws.Range("J13").Value = 0.000461987232956876
$ws.Range("M13").Value = 0.035773
$ws.Range("N13").Value = 0.107319
$ws.Range("O13").Value = 0.005010931515261538
$ws.Range("P13").Value = 0.005010931515261539
$ws.Range("Q13").Value = 0.06584845813866667
$ws.Range("R13").Value = 0.5926361232480001
$ws.Range("S13").Value = 0.000002314986385272084
$ws.Range("T13").Value = 0.000002314986385272084
$ws.Range("G14").Value = 31.11921133333334
$ws.Range("H14").Value = 93.357634
$ws.Range("I14").Value = 0.007810310653280575
$ws.Range("J14").Value = 0.007810310653280575
$ws.Range("M14").Value = 6.066605666666667
$ws.Range("N14").Value = 18.199817
$ws.Range("O14").Value = 0.849784628791665
$ws.Range("P14").Value = 0.8497846287916652
$ws.Range("Q14").Value = 188.7879838169976
$ws.Range("R14").Value = 1699.091854352978
$ws.Range("S14").Value = 0.006637081939245621
$ws.Range("T14").Value = 0.006637081939245622
$ws.Range("G15").Value = 31.11921133333334
$ws.Range("H15").Value = 93.357634
$ws.Range("I15").Value = 0.007810310653280575
$ws.Range("J15").Value = 0.007810310653280575
$ws.Range("O15").Value = 0.1196497582104962
$ws.Range("P15").Value = 0.1196497582104962
$ws.Range("Q15").Value = 26.58136644442489
$ws.Range("R15").Value = 239.232297999824
$ws.Range("S15").Value = 0.0009345017812138832
$ws.Range("T15").Value = 0.0009345017812138834
$ws.Range("G16").Value = 31.11921133333334
$ws.Range("H16").Value = 93.357634
$ws.Range("I16").Value = 0.007810310653280575
$ws.Range("J16").Value = 0.007810310653280575
$ws.Range("M16").Value = 0.1824346666666667
$ws.Range("N16").Value = 0.547304
$ws.Range("O16").Value = 0.02555468148257719
$ws.Range("P16").Value = 0.02555468148257719
$ws.Range("Q16").Value = 5.677222946526222
$ws.Range("R16").Value = 51.09500651873601
$ws.Range("S16").Value = 0.0001995900010245644
$ws.Range("T16").Value = 0.0001995900010245645
$ws.Range("G17").Value = 31.11921133333334
$ws.Range("H17").Value = 93.357634
$ws.Range("I17").Value = 0.007810310653280575
$ws.Range("J17").Value = 0.007810310653280575
$ws.Range("M17").Value = 0.035773
$ws.Range("N17").Value = 0.107319
$ws.Range("O17").Value = 0.005010931515261538
$ws.Range("P17").Value = 0.005010931515261539
$ws.Range("Q17").Value = 1.113227547027333
$ws.Range("R17").Value = 10.019047923246
$ws.Range("S17").Value = 0.00003913693179650656
$ws.Range("T17").Value = 0.00003913693179650657
$ws.Range("E18").Value = 3.0
$ws.Range("F18").Value = 1.0
$ws.Range("G18").Value = 1.905952
$ws.Range("H18").Value = 5.717856
$ws.Range("I18").Value = 0.0004783565062362683
$ws.Range("J18").Value = 0.0004783565062362683
$ws.Range("M18").Value = 6.066605666666667
$ws.Range("N18").Value = 18.199817
$ws.Range("O18").Value = 0.849784628791665
$ws.Range("P18").Value = 0.8497846287916652
$ws.Range("Q18").Value = 11.56265920359467
$ws.Range("R18").Value = 104.063932832352
$ws.Range("S18").Value = 0.0004065000060820651
$ws.Range("T18").Value = 0.0004065000060820652
$ws.Range("E19").Value = 3.0
$ws.Range("F19").Value = 1.0
$ws.Range("G19").Value = 1.905952
$ws.Range("H19").Value = 5.717856
$ws.Range("I19").Value = 0.0004783565062362683
$ws.Range("J19").Value = 0.0004783565062362683
$ws.Range("O19").Value = 0.1196497582104962
$ws.Range("P19").Value = 0.1196497582104962
$ws.Range("Q19").Value = 1.628023538090666
$ws.Range("R19").Value = 14.652211842816
$ws.Range("S19").Value = 0.00005723524030958721
$ws.Range("T19").Value = 0.00005723524030958722
$ws.Range("E20").Value = 3.0
$ws.Range("F20").Value = 1.0
$ws.Range("G20").Value = 1.905952
$ws.Range("H20").Value = 5.717856
$ws.Range("I20").Value = 0.0004783565062362683
$ws.Range("J20").Value = 0.0004783565062362683
$ws.Range("M20").Value = 0.1824346666666667
$ws.Range("N20").Value = 0.547304
$ws.Range("O20").Value = 0.02555468148257719
$ws.Range("P20").Value = 0.02555468148257719
$ws.Range("Q20").Value = 0.3477117178026667
$ws.Range("R20").Value = 3.129405460224
$ws.Range("S20").Value = 0.00001222424815198628
$ws.Range("T20").Value = 0.00001222424815198629
$ws.Range("E21").Value = 3.0
$ws.Range("F21").Value = 1.0
$ws.Range("G21").Value = 1.905952
$ws.Range("H21").Value = 5.717856
$ws.Range("I21").Value = 0.0004783565062362683
$ws.Range("J21").Value = 0.0004783565062362683
$ws.Range("M21").Value = 0.035773
$ws.Range("N21").Value = 0.107319
$ws.Range("O21").Value = 0.005010931515261538
$ws.Range("P21").Value = 0.005010931515261539
$ws.Range("Q21").Value = 0.06818162089600001
$ws.Range("R21").Value = 0.6136345880640001
$ws.Range("S21").Value = 0.00000239701169262972
$ws.Range("T21").Value = 0.00000239701169262972
